# cryptos.xlsx refresh -- Thu Jun  1 11:19:07 UTC 2023 (GitHub Actions)
#
# Re-writes the "Price" (D) and "Volume(1h)" (E) columns for every coin row
# with the latest scraped figures. Both columns are stored as plain text in
# the workbook (Price uses "." as a thousands separator, e.g. "26.874.14",
# so it is NOT a valid number; Volume is a padded "  +/-x.xx%  " string).
#
# Column D values look numeric to Excel's input parser, so a leading single
# quote is prepended to force text storage (exactly like typing an
# apostrophe before the digits into a cell) instead of letting Excel coerce
# the value into a Number and mangle it (e.g. "1.000" -> 1, "26.874.14"
# being rejected/garbled).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '26.874.14'
$ws.Range("E2").Value = '  -0.82%  '
$ws.Range("D3").Value = "'" + '1.861.66'
$ws.Range("E3").Value = '  -0.34%  '
$ws.Range("D4").Value = "'" + '1.001'
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("E5").Value = '  -0.85%  '
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("D7").Value = "'" + '0.5058'
$ws.Range("E7").Value = '  -0.70%  '
$ws.Range("D8").Value = "'" + '0.3637'
$ws.Range("E8").Value = '  -2.70%  '
$ws.Range("D9").Value = "'" + '0.07184'
$ws.Range("E9").Value = '  +0.43%  '
$ws.Range("D10").Value = "'" + '0.8957'
$ws.Range("E10").Value = '  +0.90%  '
$ws.Range("D11").Value = "'" + '20.71'
$ws.Range("E11").Value = '  +0.58%  '
$ws.Range("D12").Value = "'" + '0.07485'
$ws.Range("E12").Value = '  -0.80%  '
$ws.Range("D13").Value = "'" + '1.853.71'
$ws.Range("E13").Value = '  -0.63%  '
$ws.Range("D14").Value = "'" + '92.59'
$ws.Range("E14").Value = '  +3.67%  '
$ws.Range("E15").Value = '  -1.67%  '
$ws.Range("E16").Value = '  +0.05%  '
$ws.Range("D17").Value = "'" + '0.000008477'
$ws.Range("E17").Value = '  +0.13%  '
$ws.Range("D18").Value = "'" + '14.14'
$ws.Range("E18").Value = '  +0.15%  '
$ws.Range("E19").Value = '  -0.08%  '
$ws.Range("D20").Value = "'" + '26.912.26'
$ws.Range("E20").Value = '  -0.81%  '
$ws.Range("D21").Value = "'" + '5.029'
$ws.Range("E21").Value = '  -0.70%  '
$ws.Range("D22").Value = "'" + '2.086.99'
$ws.Range("E22").Value = '  -0.86%  '
$ws.Range("E23").Value = '  -1.98%  '
$ws.Range("D24").Value = "'" + '6.408'
$ws.Range("E24").Value = '  -1.25%  '
$ws.Range("D25").Value = "'" + '148.02'
$ws.Range("E25").Value = '  -1.68%  '
$ws.Range("D26").Value = "'" + '1.793'
$ws.Range("E26").Value = '  -2.21%  '
$ws.Range("D27").Value = "'" + '17.89'
$ws.Range("E27").Value = '  -0.21%  '
$ws.Range("D28").Value = "'" + '2.061'
$ws.Range("E28").Value = '  -1.77%  '
$ws.Range("D29").Value = "'" + '113.10'
$ws.Range("E29").Value = '  +0.49%  '
$ws.Range("D30").Value = "'" + '4.684'
$ws.Range("E30").Value = '  -1.38%  '
$ws.Range("D31").Value = "'" + '4.676'
$ws.Range("E31").Value = '  -0.20%  '
$ws.Range("D32").Value = "'" + '0.09262'
$ws.Range("E32").Value = '  +2.30%  '
$ws.Range("D33").Value = "'" + '0.05088'
$ws.Range("E33").Value = '  -0.80%  '
$ws.Range("D34").Value = "'" + '0.7440'
$ws.Range("E34").Value = '  +0.99%  '
$ws.Range("D35").Value = "'" + '2.948'
$ws.Range("E35").Value = '  -4.80%  '
$ws.Range("D36").Value = "'" + '1.149'
$ws.Range("E36").Value = '  -1.01%  '
$ws.Range("D37").Value = "'" + '3.279'
$ws.Range("E37").Value = '  +7.81%  '
$ws.Range("E38").Value = '  -1.96%  '
$ws.Range("D39").Value = "'" + '2.502'
$ws.Range("E39").Value = '  +0.68%  '
$ws.Range("D40").Value = "'" + '0.5545'
$ws.Range("E40").Value = '  +3.86%  '
$ws.Range("D41").Value = "'" + '1.069'
$ws.Range("E41").Value = '  -1.04%  '
$ws.Range("D42").Value = "'" + '118.68'
$ws.Range("E42").Value = '  +2.33%  '
$ws.Range("D43").Value = "'" + '6.484'
$ws.Range("D44").Value = "'" + '8.493'
$ws.Range("E44").Value = '  +1.74%  '
$ws.Range("E45").Value = '  -0.27%  '
$ws.Range("D46").Value = "'" + '0.4699'
$ws.Range("E46").Value = '  +1.25%  '
$ws.Range("D47").Value = "'" + '1.000'
$ws.Range("E47").Value = '  -0.09%  '
$ws.Range("D48").Value = "'" + '10.01'
$ws.Range("E48").Value = '  -0.15%  '
$ws.Range("D49").Value = "'" + '1.565'
$ws.Range("E49").Value = '  -0.02%  '
$ws.Range("D50").Value = "'" + '37.00'
$ws.Range("E50").Value = '  +1.60%  '
$ws.Range("D51").Value = "'" + '63.03'
$ws.Range("E51").Value = '  -2.34%  '
